# Aggregate discrete range test
# Adds a "Metadata" block (rows 8-12, columns A:B) to both worksheets,
# widens the relevant columns, restyles the new cells (wrap text in column B,
# Arial font for the label cells in column A), and swaps which sheet/cell is
# the active selection (sheet2 "cohort_fake" becomes the active tab).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "act score_fake"
$ws2 = $wb.Worksheets.Item(2)   # "cohort_fake"

function Add-MetadataBlock($ws) {
    $ws.Range("A8").Value = "Metadata"
    $ws.Range("B8").Value = "x"

    $ws.Range("A9").Value = "Operation-Allowed?"
    $ws.Range("B9").Value = "Yes"

    $ws.Range("A10").Value = "Sum-Allowed?"
    $ws.Range("B10").Value = "Yes"

    $ws.Range("A11").Value = "Range-Allowed?"
    $ws.Range("B11").Value = "Yes"

    $ws.Range("A12").Value = "Percentage-Allowed?"
    $ws.Range("B12").Value = "No"

    # column B of the metadata block wraps its text
    $ws.Range("B8:B12").WrapText = $true

    # column A labels (rows 9-12) use Arial instead of the default font
    $ws.Range("A9:A12").Font.Name = "Arial"
}

# Populate both sheets - order matters so the shared-string table comes out
# in the same order as the target workbook.
Add-MetadataBlock $ws1
Add-MetadataBlock $ws2

# Widen column A (labels) on both sheets, and column B on both sheets
# (values are narrower on sheet1 than sheet2).
$ws1.Columns.Item(1).ColumnWidth = 19.5
$ws1.Columns.Item(2).ColumnWidth = 11.5

$ws2.Columns.Item(1).ColumnWidth = 19.5
$ws2.Columns.Item(2).ColumnWidth = 21.67

# Move the selection on sheet1 (previously the active tab) to B10, and make
# sheet2 the active tab with C10 selected.
$ws1.Range("B10").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("C10").Select() | Out-Null

Write-Output "done"
